$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column V header, matching the existing header style (bold, bordered, centered)
$ws.Range("U1").Copy($ws.Range("V1")) | Out-Null
$ws.Range("V1").Value2 = "session_name_as_folder_name"

# Map of session titles that get a separator cleanup: "-" -> ": " (hyphen becomes colon+space)
# applied only to the *primary* separator hyphen(s), not hyphens inside compound words like
# "First-Person" or "Human-Centered".
$map = @{
    "Uncertainty Visualization-Applications, Techniques, Software, and Decision Frameworks" = "Uncertainty Visualization: Applications, Techniques, Software, and Decision Frameworks"
    "BELIV-evaluation and BEyond - methodoLogIcal approaches for Visualization" = "BELIV: evaluation and BEyond - methodoLogIcal approaches for Visualization"
    "LLM4Vis-Large Language Models for Information Visualization" = "LLM4Vis: Large Language Models for Information Visualization"
    "NLVIZ Workshop-Exploring Research Opportunities for Natural Language, Text, and Data Visualization" = "NLVIZ Workshop: Exploring Research Opportunities for Natural Language, Text, and Data Visualization"
    "EnergyVis 2024-4th Workshop on Energy Data Visualization" = "EnergyVis 2024: 4th Workshop on Energy Data Visualization"
    "VISions of the Future-Workshop on Sustainable Practices within Visualization and Physicalisation" = "VISions of the Future: Workshop on Sustainable Practices within Visualization and Physicalisation"
    "TopoInVis-Workshop on Topological Data Analysis and Visualization" = "TopoInVis: Workshop on Topological Data Analysis and Visualization"
    "Enabling Scientific Discovery-A Tutorial for Harnessing the Power of the National Science Data Fabric for Large-Scale Data Analysis" = "Enabling Scientific Discovery: A Tutorial for Harnessing the Power of the National Science Data Fabric for Large-Scale Data Analysis"
    "EduVis-2nd IEEE VIS Workshop on Visualization Education, Literacy, and Activities" = "EduVis: 2nd IEEE VIS Workshop on Visualization Education, Literacy, and Activities"
    "First-Person Visualizations for Outdoor Physical Activities-Challenges and Opportunities" = "First-Person Visualizations for Outdoor Physical Activities: Challenges and Opportunities"
    "VISxAI-7th Workshop on Visualization for AI Explainability" = "VISxAI: 7th Workshop on Visualization for AI Explainability"
    "LDAV-14th IEEE Symposium on Large Data Analysis and Visualization" = "LDAV: 14th IEEE Symposium on Large Data Analysis and Visualization"
    "VDS-Visualization in Data Science Symposium" = "VDS: Visualization in Data Science Symposium"
    "Panel-Human-Centered Computing Research in South America-Status Quo, Opportunities, and Challenges" = "Panel: Human-Centered Computing Research in South America: Status Quo, Opportunities, and Challenges"
    "Short Papers-Perception and Representation" = "Short Papers: Perception and Representation"
    "Panel-(Yet Another) Evaluation Needed? A Panel Discussion on Evaluation Trends in Visualization" = "Panel: (Yet Another) Evaluation Needed? A Panel Discussion on Evaluation Trends in Visualization"
    "Applications-Sports. Games, and Finance" = "Applications: Sports. Games, and Finance"
    "Visual Design-Sketching and Labeling" = "Visual Design: Sketching and Labeling"
    "Short Papers-Text and Multimedia" = "Short Papers: Text and Multimedia"
    "Panel-Vogue or Visionary? Current Challenges and Future Opportunities in Situated Visualizations" = "Panel: Vogue or Visionary? Current Challenges and Future Opportunities in Situated Visualizations"
    "Short Papers-Analytics and Applications" = "Short Papers: Analytics and Applications"
    "CG&A-Systems, Theory, and Evaluations" = "CG&A: Systems, Theory, and Evaluations"
    "Panel-Dear Younger Me-A Dialog About Professional Development Beyond The Initial Career Phases" = "Panel: Dear Younger Me: A Dialog About Professional Development Beyond The Initial Career Phases"
    "Applications-Industry, Computing, and Medicine" = "Applications: Industry, Computing, and Medicine"
    "Short Papers-AI and LLM" = "Short Papers: AI and LLM"
    "Application Spotlight-IEEE VIS Demos Session" = "Application Spotlight: IEEE VIS Demos Session"
    "Panel-VIS Conference Futures-Community Opinions on Recent Experiences, Challenges, and Opportunities for Hybrid Event Formats" = "Panel: VIS Conference Futures: Community Opinions on Recent Experiences, Challenges, and Opportunities for Hybrid Event Formats"
    "Panel-What Do Visualization Art Projects Bring to the VIS Community?" = "Panel: What Do Visualization Art Projects Bring to the VIS Community?"
    "Virtual-VIS from around the world" = "Virtual: VIS from around the world"
    "Short Papers-Graph, Hierarchy and Multidimensional" = "Short Papers: Graph, Hierarchy and Multidimensional"
    "Virtual-Virtual VISits" = "Virtual: Virtual VISits"
    "Panel-20 Years of Visual Analytics" = "Panel: 20 Years of Visual Analytics"
    "Short Papers-Scientific and Immersive Visualization" = "Short Papers: Scientific and Immersive Visualization"
    "CG&A-Analytics and Applications" = "CG&A: Analytics and Applications"
    "Application Spotlight-Visualization within the Department of Energy" = "Application Spotlight: Visualization within the Department of Energy"
    "Panel-Past, Present, and Future of Data Storytelling" = "Panel: Past, Present, and Future of Data Storytelling"
    "Short Papers- System design" = "Short Papers:  System design"
}

$lastRow = 324
for ($r = 2; $r -le $lastRow; $r++) {
    $old = $ws.Cells.Item($r, 1).Value2
    if ($null -eq $old) { continue }
    # New column V retains the original (pre-cleanup) session_name for folder naming
    $ws.Cells.Item($r, 22).Value2 = $old
    if ($map.ContainsKey($old)) {
        $ws.Cells.Item($r, 1).Value2 = $map[$old]
    }
}
